$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.539.23'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.566.26'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''211.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '''46.41'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.68%  '
$ws.Range("D9").Value = '''24.17'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = '1.789.75'
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("D14").Value = '1.560.96'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("D16").Value = '28.537.37'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("E17").Value = '  -3.12%  '
$ws.Range("D18").Value = '''62.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").Value = '''228.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.74%  '
$ws.Range("E20").Value = '  -2.14%  '
$ws.Range("D21").Value = '''7.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.12%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  -6.02%  '
$ws.Range("D24").Value = '''9.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.70%  '
$ws.Range("E25").Value = '  +7.60%  '
$ws.Range("D26").Value = '''150.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("E28").Value = '  -2.62%  '
$ws.Range("E29").Value = '  -3.72%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  -1.77%  '
$ws.Range("D32").Value = '''1.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").Value = '1.397.70'
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("E36").Value = '  -0.92%  '
$ws.Range("E37").Value = '  -2.97%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  +2.37%  '
$ws.Range("D40").Value = '''0.0165'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").Value = '''0.536'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("E43").Value = '  +3.11%  '
$ws.Range("D44").Value = '''0.790'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.56%  '
$ws.Range("E45").Value = '  -4.39%  '
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = '''62.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.82%  '
$ws.Range("D48").Value = '1.702.87'
$ws.Range("D49").Value = '''86.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("E50").Value = '  -4.25%  '
$ws.Range("E51").Value = '  -0.97%  '
